# Generate Report for Handoff
#
# b.md was handed off again (new .xlf files generated) for both the
# zh-cn and de-de locales. Reflect the new handoff on the Overview sheet
# and on each locale's status sheet (row 3 = b.md).

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/80ff9cdf4fbdc1cf93d4287f22288f62ce3c9593/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ca16385589a4c6e8be65b63db67c4d2a5ba22c76/e2e/b.md."

# Excel's ColumnWidth (characters) is stored internally with a constant
# +5/6 padding offset, so asking for exactly "40" character-widths stores
# as 40.8333...; back the offset out so the saved <col width="..."/> is
# exactly 40, matching the target layout.
$fortyCharWidth = 40 - (5 / 6)

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-31 04:39:18"

# ---- zh-cn sheet (b.md is row 3) ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
# "False"/"True" look like booleans to Excel's auto-typing, so they'd be
# stored as t="b" instead of the original t="s" text; use the quote-prefix
# trick to force text, then reset the style so no stray quotePrefix format
# sticks around on the cell.
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-31 04:39:13"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = $fortyCharWidth

# ---- de-de sheet (b.md is row 3) ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-31 04:39:18"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = $fortyCharWidth
